$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string for the "Total" row label - reuse the label formatting
# (Arial font) already used by the rest of column A.
$ws.Range("A43").Value = "Total"
$ws.Range("A42").Copy() | Out-Null
$ws.Range("A43").PasteSpecial(-4122) | Out-Null

# Totals row (row 43): sums over the data rows (2-42)
$ws.Range("B43").Formula = "=SUM(B2:B42)"
$ws.Range("C43").Formula = "=SUM(C2:C42)"
$ws.Range("D43").Formula = "=SUM(D2:D42)"
$ws.Range("E43").Formula = "=SUM(E2:E42)"
$ws.Range("F43").Formula = "=SUM(F2:F42)"
$ws.Range("G43").Formula = "=SUM(G2:G42)"

# Percentage row (row 44): each total column as a percentage of the grand total (B43)
$ws.Range("C44").Formula = "=C43/B43*100"
$ws.Range("D44").Formula = "=D43/B43*100"
$ws.Range("E44").Formula = "=E43/B43*100"
$ws.Range("F44").Formula = "=F43/B43*100"
$ws.Range("G44").Formula = "=G43/B43*100"

# Update the used dimension / view to reflect the newly added rows - scroll
# the window so row 39 is the top visible row, then select F54.
$ws.Range("A39").Select()
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F54").Select()
